$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.093.69"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "3.806.33"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'708.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +12.46%  "
$ws.Range("D6").Value = "'174.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("D7").Value = "3.804.55"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("D11").Value = "'7.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.79%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("E13").Value = "  +7.84%  "
$ws.Range("D14").Value = "'36.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.80%  "
$ws.Range("D15").Value = "4.443.50"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "3.804.81"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "71.068.86"
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'11.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +18.75%  "
$ws.Range("D22").Value = "'484.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("D23").Value = "'0.719"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "'83.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.30%  "
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "'12.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.01%  "
$ws.Range("D27").Value = "'10.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.78%  "
$ws.Range("D28").Value = "'2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.54%  "
$ws.Range("D29").Value = "3.954.17"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'3.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.60%  "
$ws.Range("D32").Value = "'7.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.81%  "
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("E34").Value = "  +4.82%  "
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("D36").Value = "'9.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.46%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "3.754.54"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("E39").Value = "  +3.48%  "
$ws.Range("D40").Value = "'3.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.62%  "
$ws.Range("D41").Value = "'6.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.42%  "
$ws.Range("E42").Value = "  +14.96%  "
$ws.Range("D43").Value = "'0.000333"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +26.27%  "
$ws.Range("D44").Value = "'0.971"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D47").Value = "'45.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.24%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'49.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.48%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'160.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").Value = "'1.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("E51").Value = "  +2.89%  "
